# Append new asset-tracking events to the "4.2 Timestamps" log and the
# "All SANs" summary sheet.
#
#   4.2 Timestamps!A69:D69 -> 2024-01-10 19:22:14 | Laptop 840 G9   | add | SAN123456
#   4.2 Timestamps!A70:D70 -> 2024-01-11 20:15:00 | Desktop Mini G9 | add | SAN123432
#   All SANs!A5:C5         -> SAN123456 | Laptop 840 G9   | 2024-01-10 19:22:14
#   All SANs!A6:C6         -> SAN123432 | Desktop Mini G9 | 2024-01-11 20:15:00
#
# Also re-stamps rows 67/68 of "4.2 Timestamps" (the two most recent,
# pre-existing log entries) with the same formatting used by every other
# data row in that table, so the whole tail of the log is visually
# consistent.

$wb = $excel.ActiveWorkbook

$xlGeneral = 1
$xlBottom = -4107

# ---------------------------------------------------------------------
# "4.2 Timestamps" - the running add/subtract log.
# ---------------------------------------------------------------------
$tsLog = $wb.Worksheets.Item("4.2 Timestamps")

# Re-enter the last two existing rows so they pick up the same look as
# the rest of the table (this also guarantees the whole block 67:70 is
# formatted uniformly before the two brand-new rows are appended below).
$tsLog.Rows("67:68").Delete()

$tsLog.Range("A67").Value = "2024-01-03 21:44:59"
$tsLog.Range("B67").Value = "Desktop Mini G9"
$tsLog.Range("C67").Value = "add"
$tsLog.Range("D67").Value = "SAN147896"

$tsLog.Range("A68").Value = "2024-01-03 21:45:07"
$tsLog.Range("B68").Value = "Desktop Mini G9"
$tsLog.Range("C68").Value = "subtract"
$tsLog.Range("D68").Value = "SAN147896"

# New event: Laptop 840 G9 / SAN123456 added.
$tsLog.Range("A69").Value = "2024-01-10 19:22:14"
$tsLog.Range("B69").Value = "Laptop 840 G9"
$tsLog.Range("C69").Value = "add"
$tsLog.Range("D69").Value = "SAN123456"

$tsLog.Range("A67:D69").HorizontalAlignment = $xlGeneral
$tsLog.Range("A67:D69").VerticalAlignment = $xlBottom

# New event: Desktop Mini G9 / SAN123432 added - left in the sheet's
# plain/default look, matching how it was entered.
$tsLog.Range("A70").Value = "2024-01-11 20:15:00"
$tsLog.Range("B70").Value = "Desktop Mini G9"
$tsLog.Range("C70").Value = "add"
$tsLog.Range("D70").Value = "SAN123432"

$tsLog.Range("A70:D70").HorizontalAlignment = $xlGeneral
$tsLog.Range("A70:D70").VerticalAlignment = $xlBottom

# ---------------------------------------------------------------------
# "All SANs" - flat summary of every SAN ever logged.
# ---------------------------------------------------------------------
$allSans = $wb.Worksheets.Item("All SANs")

$allSans.Range("A5").Value = "SAN123456"
$allSans.Range("B5").Value = "Laptop 840 G9"
$allSans.Range("C5").Value = "2024-01-10 19:22:14"

$allSans.Range("A6").Value = "SAN123432"
$allSans.Range("B6").Value = "Desktop Mini G9"
$allSans.Range("C6").Value = "2024-01-11 20:15:00"

$allSans.Range("A5:C6").HorizontalAlignment = $xlGeneral
$allSans.Range("A5:C6").VerticalAlignment = $xlBottom
